# Update the Hypothesis statement on the "Add your Hypotheses" slide:
# spell out the abbreviated variable names PS/G and MP in both the
# Null-hypothesis and Alternative-hypothesis example sentences.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# Null hypothesis example sentence (paragraph 4):
#   "... no correlation between PS/G and MP."
$tr.Runs(16, 1).Text = "Point Scored/Game"
$tr.Runs(18, 1).Text = "Minutes Played"

# Alternative hypothesis example sentence (paragraph 6):
#   "... a correlation between PS/G and MP."
$tr.Runs(29, 1).Text = "Point Scored/Game"
$tr.Runs(31, 1).Text = "Minutes Played"
